# TC24_CDS_Filter_PHSAccession-phs002517_LibLyot_InsModel_RefGenAsmbly.xlsx
#
# The underlying "startup" sheet content (tab names, queries, labels) is not
# changing - the workbook was simply re-opened/re-saved (picking up a newer
# Office build's auto-generated version stamps, which this headless runtime
# regenerates on its own and which aren't reachable through the Excel object
# model), and the user's last selection before saving was cell C3 on the
# "startup" sheet. Reproduce that by moving the active selection there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$ws.Activate()
$ws.Range("C3").Select()
